$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update staff record ---
# Admin A001 had its password (salt + hash) updated.
$ws.Range("B8").Value = "ea6697e927d5087d08b5218455deb794"
$ws.Range("C8").Value = "a74c91fe1367f03956a8afb4bdd0c201982122e1cf3fdd09ea2445e8cccaf5a6"

# --- Add staff records ---
# New doctor D999.
$ws.Range("A9").Value = "D999"
$ws.Range("B9").Value = "e809a49832b33d35c47689286b1cfbaa"
$ws.Range("C9").Value = "85da2685c4bb4f84dfcb126a3ca4bc958ab9c8c27b24c80485fb34731f64a71e"

# New patient P999.
$ws.Range("A10").Value = "P999"
$ws.Range("B10").Value = "f86bf1e8753ab0a04d292050e6841f0e"
$ws.Range("C10").Value = "b2597e074e2ab8822cc65737403429b1f1fa8139313b0907e05c316816bba347"

# Reflect the two newly-added rows in the sheet's current selection.
$ws.Range("A9:XFD10").Select()
